# Weekly update: insert a new price record for "Locoto" (Vega Modelo de
# Temuco) at the top of the data block (new row 18), pushing the existing
# rows 18-34 down to 19-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18; Excel copies formatting
# (including the date cell's number format in column D) down from the
# row that used to be 18 (now 19), matching native Excel behaviour.
$ws.Range("A18:R18").EntireRow.Insert()

# Populate the new row 18 with this week's record.
$ws.Range("A18").Value2 = 10
$ws.Range("B18").Value2 = "Vega Modelo de Temuco"
$ws.Range("C18").Value2 = "La Araucanía"
$ws.Range("D18").Value2 = 44781
$ws.Range("E18").Value2 = 9
$ws.Range("F18").Value2 = 100112042
$ws.Range("G18").Value2 = "Locoto"
$ws.Range("H18").Value2 = "Sin especificar"
$ws.Range("I18").Value2 = "Primera"
$ws.Range("J18").Value2 = 250
$ws.Range("K18").Value2 = 2700
$ws.Range("L18").Value2 = 2700
$ws.Range("M18").Value2 = 2700
$ws.Range("N18").Value2 = "$/kilo"
$ws.Range("O18").Value2 = "Región de Arica y Parinacota"
$ws.Range("P18").Value2 = 2700
$ws.Range("Q18").Value2 = 1
$ws.Range("R18").Value2 = "Hortaliza"
